$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "mec": add rows for m4, m5, m6 (each with value 1)
# ---------------------------------------------------------------
$wsMec = $wb.Worksheets.Item("mec")
$wsMec.Cells.Item(5, 1).Value = "m4"
$wsMec.Cells.Item(5, 2).Value = 1
$wsMec.Cells.Item(6, 1).Value = "m5"
$wsMec.Cells.Item(6, 2).Value = 1
$wsMec.Cells.Item(7, 1).Value = "m6"
$wsMec.Cells.Item(7, 2).Value = 1
[void]$wsMec.Range("A8").Select()

# ---------------------------------------------------------------
# Sheet "service": scale up existing B/C values for rows 2-7
# ---------------------------------------------------------------
$wsService = $wb.Worksheets.Item("service")
$wsService.Cells.Item(2, 2).Value = 180
$wsService.Cells.Item(2, 3).Value = 9000
$wsService.Cells.Item(3, 2).Value = 44
$wsService.Cells.Item(3, 3).Value = 9000
$wsService.Cells.Item(4, 2).Value = 20
$wsService.Cells.Item(4, 3).Value = 9000
$wsService.Cells.Item(5, 2).Value = 180
$wsService.Cells.Item(5, 3).Value = 4000
$wsService.Cells.Item(6, 2).Value = 44
$wsService.Cells.Item(6, 3).Value = 4000
$wsService.Cells.Item(7, 2).Value = 20
$wsService.Cells.Item(7, 3).Value = 4000
[void]$wsService.Range("B8").Select()

# ---------------------------------------------------------------
# Sheet "request": update values for rows 2-4 and add rows 5-7 (m4, m5, m6)
# ---------------------------------------------------------------
$wsRequest = $wb.Worksheets.Item("request")

$wsRequest.Cells.Item(2, 2).Value = 250
$wsRequest.Cells.Item(2, 3).Value = 300
$wsRequest.Cells.Item(2, 4).Value = 700
$wsRequest.Cells.Item(2, 5).Value = 250
$wsRequest.Cells.Item(2, 6).Value = 300
$wsRequest.Cells.Item(2, 7).Value = 700

$wsRequest.Cells.Item(3, 2).Value = 100
$wsRequest.Cells.Item(3, 3).Value = 300
$wsRequest.Cells.Item(3, 4).Value = 1000
$wsRequest.Cells.Item(3, 5).Value = 100
$wsRequest.Cells.Item(3, 6).Value = 300
$wsRequest.Cells.Item(3, 7).Value = 1000

$wsRequest.Cells.Item(4, 2).Value = 150
$wsRequest.Cells.Item(4, 3).Value = 500
$wsRequest.Cells.Item(4, 4).Value = 700
$wsRequest.Cells.Item(4, 5).Value = 150
$wsRequest.Cells.Item(4, 6).Value = 500
$wsRequest.Cells.Item(4, 7).Value = 700

$wsRequest.Cells.Item(5, 1).Value = "m4"
$wsRequest.Cells.Item(5, 2).Value = 100
$wsRequest.Cells.Item(5, 3).Value = 100
$wsRequest.Cells.Item(5, 4).Value = 100
$wsRequest.Cells.Item(5, 5).Value = 100
$wsRequest.Cells.Item(5, 6).Value = 100
$wsRequest.Cells.Item(5, 7).Value = 100

$wsRequest.Cells.Item(6, 1).Value = "m5"
$wsRequest.Cells.Item(6, 2).Value = 100
$wsRequest.Cells.Item(6, 3).Value = 100
$wsRequest.Cells.Item(6, 4).Value = 100
$wsRequest.Cells.Item(6, 5).Value = 100
$wsRequest.Cells.Item(6, 6).Value = 100
$wsRequest.Cells.Item(6, 7).Value = 100

$wsRequest.Cells.Item(7, 1).Value = "m6"
$wsRequest.Cells.Item(7, 2).Value = 100
$wsRequest.Cells.Item(7, 3).Value = 100
$wsRequest.Cells.Item(7, 4).Value = 100
$wsRequest.Cells.Item(7, 5).Value = 100
$wsRequest.Cells.Item(7, 6).Value = 100
$wsRequest.Cells.Item(7, 7).Value = 100
[void]$wsRequest.Range("A8").Select()

# ---------------------------------------------------------------
# Sheet "intervalForSendingRequests": no value changes needed
# (cell content unchanged; shared-string index shift only)
# ---------------------------------------------------------------
$wsInterval = $wb.Worksheets.Item("intervalForSendingRequests")
[void]$wsInterval.Range("D7").Select()

# ---------------------------------------------------------------
# Restore "request" as the active sheet (it was the selected tab
# before these edits and remains so afterward).
# ---------------------------------------------------------------
$wsRequest.Activate()
